$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 0.05338078291814947
$ws.Range("C2").Value = 0.05
$ws.Range("E2").Value = 0.09523809523809523
$ws.Range("F2").Value = 0.2083333333333333
$ws.Range("G2").Value = 0.5777777777777777
$ws.Range("H2").Value = 0.6830524344569288
$ws.Range("J2").Value = 532
$ws.Range("K2").Value = 2

# --- Classification Report sheet ---
$ws = $wb.Worksheets.Item("Classification Report")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.003745318352059925
$ws.Range("D2").Value = 0.007462686567164179

$ws.Range("B3").Value = 0.05
$ws.Range("D3").Value = 0.09523809523809523

$ws.Range("B4").Value = 0.05338078291814947
$ws.Range("C4").Value = 0.05338078291814947
$ws.Range("D4").Value = 0.05338078291814947
$ws.Range("E4").Value = 0.05338078291814947

$ws.Range("B5").Value = 0.525
$ws.Range("C5").Value = 0.50187265917603
$ws.Range("D5").Value = 0.0513503909026297

$ws.Range("B6").Value = 0.9526690391459074
$ws.Range("C6").Value = 0.05338078291814947
$ws.Range("D6").Value = 0.01183583860059135

# --- Confusion Matrix sheet ---
$ws = $wb.Worksheets.Item("Confusion Matrix")
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 532
